# Add I0 and IF columns (I and J) to the worksheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for columns I and J, matching the style used by existing headers (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:J44
$iValues = @(2,7,4,7,8,7,6,6,7,7,7,5,4,7,5,5,6,5,8,2,9,5,5,6,1,1,7,8,7,5,6,6,8,6,6,6,3,8,6,6,8,8,3)
$jValues = @(3,8,5,7,8,8,6,7,8,8,7,5,5,8,6,6,7,6,8,5,9,7,5,6,3,3,8,8,7,6,7,7,9,6,6,9,6,9,7,7,8,9,4)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
